$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Test Suite object result for C2 and C3 failed -> change from "Y" to "N"
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"

# Update the active selection to C2
$ws.Range("C2").Select()
